# Update omega-3 row: correct quantity/scale.
# The unit label changes from mg to g, and the threshold value is corrected
# from 250 (mg) to 2.5 (g) expressed directly in grams (no /1000 conversion
# formula needed anymore since the source value is already in grams).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "o3(g)"
$ws.Range("D8").Value = 2.5
$ws.Range("E8").Value = 2.5

# Reflect the author's final cell selection (cursor left on E9 after editing).
$ws.Range("E9").Select() | Out-Null
